{"js": "// The Jekyll-generated site footer that trails the \"M\u00c1QUINASEscola PRO-TEC\"\n// paragraph is removed on rebuild:\n//   - a blank paragraph\n//   - \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   - \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github\n//     pages. Original theme under Creative Commons Attribution\"\n// The paragraph that follows them (another blank paragraph right before the\n// page-break paragraph) is left untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targets = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\n// Locate the \"Ver no Jupiter...\" paragraph; the blank paragraph immediately\n// before it and the copyright paragraph immediately after it go with it.\nconst items = paragraphs.items;\nlet startIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === targets[0]) {\n    startIndex = i;\n    break;\n  }\n}\n\nif (\n  startIndex > 0 &&\n  items[startIndex - 1] && items[startIndex - 1].text === \"\" &&\n  items[startIndex + 1] && items[startIndex + 1].text === targets[1]\n) {\n  // Delete from the end backwards so earlier indices stay valid.\n  items[startIndex + 1].delete(); // \"\u00a9 2020 ...\" paragraph\n  items[startIndex].delete();     // \"Ver no Jupiter ...\" paragraph\n  items[startIndex - 1].delete(); // preceding blank paragraph\n  await context.sync();\n}\n", "ps1": "# The Jekyll-generated site footer that trails the \"MAQUINASEscola PRO-TEC\"\n# paragraph is removed on rebuild:\n#   - a blank paragraph\n#   - \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   - \"(c) 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github\n#     pages. Original theme under Creative Commons Attribution\"\n# The paragraph that follows them (another blank paragraph right before the\n# page-break paragraph) is left untouched.\n\n$d = $word.ActiveDocument\n\n$marker = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n\n$targetIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $text = $d.Paragraphs.Item($i).Range.Text.TrimEnd()\n    if ($text -eq $marker) {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -gt 1 -and $targetIndex -lt $count) {\n    $prevText = $d.Paragraphs.Item($targetIndex - 1).Range.Text.TrimEnd()\n    $nextText = $d.Paragraphs.Item($targetIndex + 1).Range.Text.TrimEnd()\n\n    if ($prevText -eq \"\" -and $nextText.Contains(\"luizeleno@usp.br\")) {\n        # Delete the three paragraphs from the end backwards so earlier\n        # indices remain valid while deleting.\n        $d.Paragraphs.Item($targetIndex + 1).Range.Delete() # copyright / Jekyll footer line\n        $d.Paragraphs.Item($targetIndex).Range.Delete()     # \"Ver no Jupiter...\" line\n        $d.Paragraphs.Item($targetIndex - 1).Range.Delete() # preceding blank paragraph\n    }\n}\n"}
